$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the column headers in row 1 -----------------------------
# "_old" -> "_FV2410", "_new" -> "_FV2504" (the "diff" column, K1, is unchanged)
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val.ToString().Replace("_old", "_FV2410").Replace("_new", "_FV2504")
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# --- 2) Turn the used range into an Excel Table (ListObject) -----------
$rng = $ws.Range("A1:U75")
$tbl = $ws.ListObjects.Add(1, $rng, $false, 1)
$tbl.Name = "Table1"
# Match the source file's plain (unnamed) table style - no banding color.
$tbl.TableStyle = ""

# --- 3) Freeze the header row -------------------------------------------
# Selecting the cell below/right of the freeze point and then toggling
# FreezePanes is how Excel itself records a plain top-row freeze
# (as opposed to explicitly driving SplitRow/SplitColumn, which yields a
# "frozenSplit" pane instead of "frozen").
$ws.Range("A2").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
